# Generate Report for Archive
# Localization status moved from "Ready for handoff" to "In Translation".
# Update the status text everywhere it appears (Overview + per-locale sheets)
# and re-fit the status columns, since their content length changed.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # NOTE: keep the known literal on the LEFT of -eq. Value2 can come
            # back as a native Boolean for boolean-looking cells (e.g. "True"),
            # and PowerShell's -eq coerces the RHS to the LHS's type, so
            # "$cell.Value2 -eq $oldStatus" would wrongly coerce $oldStatus to
            # $true and match every boolean cell.
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
                $cell.EntireColumn.AutoFit() | Out-Null
            }
        }
    }
}
